# Applies updated NOAA temperature data (column K: average_county_temperature)
# and the resulting recalculated worst/best ASHP COP values (columns R and S)
# for facility rows 2, 6, 7 and 10 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (facility_id 1000032)
$ws.Range("K2").Value = 21.28240740740739
$ws.Range("R2").Value = 2.017497406510892
$ws.Range("S2").Value = 2.228623569098047

# Row 6 (facility_id 1006856)
$ws.Range("K6").Value = 19.65277777777778
$ws.Range("R6").Value = 1.983015294974508
$ws.Range("S6").Value = 2.18606997558991

# Row 7 (facility_id 1006919)
$ws.Range("K7").Value = 13.75752314814816
$ws.Range("R7").Value = 1.867546171126113
$ws.Range("S7").Value = 2.044826120875009

# Row 10 (facility_id 1013817)
$ws.Range("K10").Value = 19.65277777777778
$ws.Range("R10").Value = 1.983015294974508
$ws.Range("S10").Value = 2.18606997558991
